# Retirer la décote si l'impôt est déjà nul
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reset the three income inputs to 0
$ws.Range("A2").Value = 0
$ws.Range("A4").Value = 0
$ws.Range("A6").Value = 0

# Switch household status from "Couple" to "Célibataire" (and the related
# parts count used elsewhere in the sheet)
$ws.Range("A8").Value = "Célibataire"
$ws.Range("D8").Value = 1

# Only apply the "décote" (tax rebate) when the computed tax (row 17) is
# strictly positive; otherwise there is nothing to rebate.
$ws.Range("E18").Formula = '=IF(E17>0,IF($A$8="Couple",IF(E17 <$A$20,($A$18 - ($D$18 * E17)),0),IF(E17 <$A$15,($A$13 - ($D$18 * E17)),0)),0)'
$ws.Range("F18").Formula = '=IF(F17>0,IF($A$8="Couple",IF(F17 <$A$20,($A$18 - ($D$18 * F17)),0),IF(F17 <$A$15,($A$13 - ($D$18 * F17)),0)),0)'
$ws.Range("G18").Formula = '=IF(G17>0,IF($A$8="Couple",IF(G17 <$A$20,($A$18 - ($D$18 * G17)),0),IF(G17 <$A$15,($A$13 - ($D$18 * G17)),0)),0)'
$ws.Range("H18").Formula = '=IF(H17>0,IF($A$8="Couple",IF(H17 <$A$20,($A$18 - ($D$18 * H17)),0),IF(H17 <$A$15,($A$13 - ($D$18 * H17)),0)),0)'

# Restore the original cell selection location
$ws.Range("A7").Select()

$wb.Save()
